# TFS 3027: Adding rows for SrMgr Dashboard support
# (Role, Entitlement, Role_Entitlement Link, Role Access Tables)

$wb = $excel.ActiveWorkbook

# --- AT_Role: add RoleId 105 "SeniorManager" ---
$wsRole = $wb.Worksheets.Item("AT_Role")
$wsRole.Range("A6").Value = 105
$wsRole.Range("B6").Value = "SeniorManager"
$wsRole.Range("C6").Value = 0
$wsRole.Range("A2").Select()

# --- AT_Entitlement: add EntitlementId 208 "SeniorManagerDashboard" ---
$wsEnt = $wb.Worksheets.Item("AT_Entitlement")
$wsEnt.Range("A9").Value = 208
$wsEnt.Range("B9").Value = "SeniorManagerDashboard"
$wsEnt.Range("A2").Select()

# --- AT_Role_Entitlement_Link: link RoleId 105 to EntitlementId 208 ---
$wsLink = $wb.Worksheets.Item("AT_Role_Entitlement_Link")
$wsLink.Range("A14").Value = 105
$wsLink.Range("B14").Value = 208
$wsLink.Range("A2").Select()

# --- AT_Role_Access: re-sort existing rows and append new SeniorManager row ---
$wsAccess = $wb.Worksheets.Item("AT_Role_Access")

$wsAccess.Range("A2").Value = "WISY14"
$wsAccess.Range("B2").Value = "Principal Analyst, Systems"
$wsAccess.Range("C2").Value = 101
$wsAccess.Range("D2").Value = "CoachingAdmin"
$wsAccess.Range("E2").Value = 0
$wsAccess.Range("F2").Value = 1

$wsAccess.Range("A3").Value = "WACQ13"
$wsAccess.Range("B3").Value = "Sr Specialist, Quality (CS)"
$wsAccess.Range("C3").Value = 101
$wsAccess.Range("D3").Value = "CoachingAdmin"
$wsAccess.Range("E3").Value = 0
$wsAccess.Range("F3").Value = 1

$wsAccess.Range("A4").Value = "WACS50"
$wsAccess.Range("B4").Value = "Manager, Customer Service"
$wsAccess.Range("C4").Value = 102
$wsAccess.Range("D4").Value = "CoachingUser"
$wsAccess.Range("E4").Value = 1
$wsAccess.Range("F4").Value = 1

$wsAccess.Range("A5").Value = "WACS60"
$wsAccess.Range("B5").Value = "Sr Manager, Customer Service"
$wsAccess.Range("C5").Value = 102
$wsAccess.Range("D5").Value = "CoachingUser"
$wsAccess.Range("E5").Value = 1
$wsAccess.Range("F5").Value = 1

$wsAccess.Range("A6").Value = "WIHD50"
$wsAccess.Range("B6").Value = "Manager, Help Desk"
$wsAccess.Range("C6").Value = 102
$wsAccess.Range("D6").Value = "CoachingUser"
$wsAccess.Range("E6").Value = 1
$wsAccess.Range("F6").Value = 1

$wsAccess.Range("A7").Value = "WTTR50"
$wsAccess.Range("B7").Value = "Manager, Training"
$wsAccess.Range("C7").Value = 102
$wsAccess.Range("D7").Value = "CoachingUser"
$wsAccess.Range("E7").Value = 1
$wsAccess.Range("F7").Value = 1

$wsAccess.Range("A8").Value = "WPPM13"
$wsAccess.Range("B8").Value = "Sr Analyst, Program"
$wsAccess.Range("C8").Value = 102
$wsAccess.Range("D8").Value = "CoachingUser"
$wsAccess.Range("E8").Value = 1
$wsAccess.Range("F8").Value = 1

$wsAccess.Range("A9").Value = "WISY14"
$wsAccess.Range("B9").Value = "Principal Analyst, Systems"
$wsAccess.Range("C9").Value = 103
$wsAccess.Range("D9").Value = "WarningAdmin"
$wsAccess.Range("E9").Value = 0
$wsAccess.Range("F9").Value = 1

$wsAccess.Range("A10").Value = "WACQ13"
$wsAccess.Range("B10").Value = "Sr Specialist, Quality (CS)"
$wsAccess.Range("C10").Value = 103
$wsAccess.Range("D10").Value = "WarningAdmin"
$wsAccess.Range("E10").Value = 0
$wsAccess.Range("F10").Value = 1

$wsAccess.Range("A11").Value = "WACS60"
$wsAccess.Range("B11").Value = "Sr Manager, Customer Service"
$wsAccess.Range("C11").Value = 105
$wsAccess.Range("D11").Value = "SeniorManager"
$wsAccess.Range("E11").Value = 1
$wsAccess.Range("F11").Value = 1

$wsAccess.Range("A2").Select()

# --- Revision_History: log this change as revision 4, TFS 3027 ---
$wsRev = $wb.Worksheets.Item("Revision_History")
$wsRev.Range("A5").Copy($wsRev.Range("A6"))
$wsRev.Range("A6").Value = 4
$wsRev.Range("B5").Copy($wsRev.Range("B6"))
$wsRev.Range("B6").Value = 42691
$wsRev.Range("C6").Value = "Susmitha Palacherla"
$wsRev.Range("D6").Value = 3027
$wsRev.Range("E6").Value = "Adding rows for SrMgr Dashboard support(Role, Entitlement, Role_Entitlement Link, Role Access Tables"
$wsRev.Range("E6").Select()
